$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2, A3, A4 with the newly merged/concatenated text values
$ws.Range("A2").Value = "(`"City's Blessing`", ['Card', 'Elemental', 'Token Creature — Elemental', 'At the beginning of your upkeep, sacrifice this creature and return target card named Rekindling Phoenix from your graveyard to the battlefield. It gains haste until end of turn.', '0/1'])"

$ws.Range("A3").Value = "('Illusion', ['Token Creature — Illusion', 'When this creature becomes the target of a spell, sacrifice it.', '2/2', 'Saproling', 'Token Creature — Saproling', '1/1'])"

$ws.Range("A4").Value = "('Merfolk', ['Token Creature — Merfolk', 'Hexproof (This creature can’t be the target of spells or abilities your opponents control.)', '1/1', 'Treasure', 'Token Artifact — Treasure', '{T}, Sacrifice this artifact: Add one mana of any color.'])"

# Remove the now-obsolete rows 5 through 21 (their content moved into A2:A4 above)
$ws.Range("A5:A21").ClearContents()
